$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# Note: the engine stores column width using Excel's internal
# "characters -> pixels -> characters" round-trip conversion, which adds
# a small padding offset. The input values below are chosen so the
# resulting stored OOXML width matches the target (30 / 66) exactly.
$ws.Columns.Item(1).ColumnWidth = 29.16
$ws.Columns.Item(2).ColumnWidth = 65.16

# Update cell text values
$ws.Range("A1").Value = "button_showApiKey_trNthChild"
$ws.Range("B2").Value = "Data Files/AI-Generated/Common/createNewApiKeyWithName-test-data"
